# Commit: "Doing Updates for Financials"
#
# The WSFS yearly-financials sheet gets a new reporting-period column
# inserted right after "Period Ending" / the label column (i.e. a new
# column D), pushing the previously-existing D:K figures out to E:L.
# The freshly inserted column D is then filled in with that period's
# Income Statement / Balance Sheet / Cash Flow figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WSFS")

# Insert a new column at D; the former D:K data shifts right to E:L.
$ws.Columns.Item(4).Insert()

# The new column D should carry the same number formatting as its
# neighbor column E (dates on the "Period Ending" rows, plain numbers
# everywhere else), so copy formats across before writing values.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# A freshly inserted column otherwise falls back to the sheet's default
# width; match it back up with the rest of the data columns.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth()

# Populate the new column D with the latest reporting period's data.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 293000
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 59700
$ws.Range("D18").Value = 233300
$ws.Range("D20").Value = -62500
$ws.Range("D21").Value = 182100
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 170800
$ws.Range("D24").Value = 36600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 134200
$ws.Range("D27").Value = 134200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 500
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 62500
$ws.Range("D33").Value = 134700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 134700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 619600
$ws.Range("D42").Value = 57700
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 45000
$ws.Range("D49").Value = 186000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7248900
$ws.Range("D57").Value = 1900
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 213300
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 6427900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 791000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 820900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 134700
$ws.Range("D83").Value = 11300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 135600
$ws.Range("D91").Value = -5500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -359900
$ws.Range("D96").Value = -13200
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 121200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -103100
